# The commit removes the record for Caso 6029 (TANDIL 4746, row 15) from the
# "PEBCOM" sheet. All subsequent rows shift up by one, and the sheet's used
# range shrinks from A1:P49 to A1:P48.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 15 (Caso 6029 / TANDIL 4746) and shift the rows below
# it upward, just like pressing Ctrl+"-" on a selected row in Excel.
$ws.Rows.Item(15).Delete()
